$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("task")
$ws.Activate()

# Insert a new column before column Q (17) - shifts Q..S -> R..T
$ws.Columns("Q").Insert()

# New column header (row 1) and new item reward entry (row 6)
$ws.Range("Q1").Value = "Item`nitemid,quantity,당첨범위,전체범위"
$ws.Range("Q6").Value = "item-1,3,10,10`nitem-2,1,50,100"
$ws.Range("Q6").WrapText = $true

# Match the authored column width for the new column
$ws.Columns("Q").ColumnWidth = 22

$ws.Range("Q6").Select() | Out-Null
